# Applies updated blob-detection values to Sheet1, per commit:
# "Will erase most blob detection settings to minimize the error"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 5.077576160430908
$ws.Range("C4").Value  = 3.634586334228516
$ws.Range("C5").Value  = 3.613167524337769
$ws.Range("C6").Value  = 2.955438852310181
$ws.Range("C9").Value  = 7.289011478424072
$ws.Range("C10").Value = 2.955438852310181
$ws.Range("C11").Value = 9.361291885375977
$ws.Range("C12").Value = 10.91279149055481
$ws.Range("C13").Value = 16.16635513305664
$ws.Range("C14").Value = 8.831316709518433
$ws.Range("C16").Value = 4.448910713195801
$ws.Range("C17").Value = 11.33061385154724
$ws.Range("C18").Value = 2.955438852310181
$ws.Range("C20").Value = 7.917736530303955
$ws.Range("C21").Value = 14.84580516815186
$ws.Range("C22").Value = 22.16390085220337
$ws.Range("C24").Value = 12.40105843544006
$ws.Range("C25").Value = 8.337146282196045
$ws.Range("C26").Value = 4.315976142883301
$ws.Range("C28").Value = 11.40879583358765
$ws.Range("D28").Value = 0
$ws.Range("C32").Value = 31.79722595214844
